$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EQUATION column (B/G): Helmholtz -> Burgers
$ws.Range("B2").Value = "Burgers"
$ws.Range("G2").Value = "Burgers"

# MODEL column (C): ResNet -> NetA
$ws.Range("C2").Value = "NetA"

# FOLDER column (F): 100N31 -> 1000N31
$ws.Range("F2").Value = "1000N31"

# TIMESTAMP (E2)
$ws.Range("E2").Value = 44051.19191161724

# BLOCKS (I2)
$ws.Range("I2").Value = 2

# BATCH (L2)
$ws.Range("L2").Value = 1000

# EPOCHS (M2)
$ws.Range("M2").Value = 100000

# AVG IT/S (N2)
$ws.Range("N2").Value = 0.3

# LOSS (O2)
$ws.Range("O2").Value = 2.304236

# MAEa (P2)
$ws.Range("P2").Value = 0.0008630000000000001

# MSEa (Q2)
$ws.Range("Q2").Value = 0.050889

# MIEa (R2)
$ws.Range("R2").Value = 0.045676

# MAEu (S2)
$ws.Range("S2").Value = 0.004426

# MSEu (T2)
$ws.Range("T2").Value = 0.044978

# MIEu (U2)
$ws.Range("U2").Value = 0.045374

# NBFUNCS (V2)
$ws.Range("V2").Value = 2
